$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.518.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.07%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.825.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.06%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'600.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.33%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'163.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.00%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.822.56"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +2.05%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.02%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -2.42%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -3.51%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -1.11%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -0.74%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'36.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -3.87%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -2.54%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.460.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.01%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.827.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +2.13%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'68.633.21"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.85%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'7.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.69%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -0.55%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'17.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.94%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'11.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.51%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'484.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.82%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -1.70%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.0000160"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +6.62%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'83.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.11%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.71%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'12.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.98%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'9.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.82%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -0.29%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -1.44%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'7.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -4.47%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.972.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.04%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -4.20%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'31.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.33%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'3.766.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E36").Value = "'  -1.50%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +0.86%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -0.80%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -2.31%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.05%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -3.23%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'434.61"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.49%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.95"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -4.15%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'48.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.79%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -0.95%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -0.02%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'8.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.91%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.835.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.40%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'142.62"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.97%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +0.16%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +13.30%  "
$ws.Range("E51").Style = "Normal"
